# Atualização de bases das ligas, do dia: 17-03-2024 às 10:24
# Swap the per-match data (columns B:AC) between each pair of rows listed
# below, while leaving column A (the sequence/id column) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(22, 23),
    @(54, 55),
    @(58, 59),
    @(73, 74),
    @(78, 79),
    @(91, 92),
    @(108, 109)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B${r1}:AC${r1}")
    $range2 = $ws.Range("B${r2}:AC${r2}")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
